# Insert a new data row before row 333 (shifting existing rows 333..431 down
# to 334..432) and populate it with the new "Papa" price record described in
# the commit. This mirrors Excel's Rows.Insert with xlShiftDown (-4121).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(333).Insert(-4121)

$ws.Cells.Item(333, 1).Value  = 8
$ws.Cells.Item(333, 2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(333, 3).Value  = "Coquimbo"
$ws.Cells.Item(333, 4).Value  = 44736
$ws.Cells.Item(333, 5).Value  = 4
$ws.Cells.Item(333, 6).Value  = 100114001
$ws.Cells.Item(333, 7).Value  = "Papa"
$ws.Cells.Item(333, 8).Value  = "Cardinal"
$ws.Cells.Item(333, 9).Value  = "1a nueva(o)"
$ws.Cells.Item(333, 10).Value = 2400
$ws.Cells.Item(333, 11).Value = 9000
$ws.Cells.Item(333, 12).Value = 10000
$ws.Cells.Item(333, 13).Value = 9500
$ws.Cells.Item(333, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(333, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(333, 16).Value = 380
$ws.Cells.Item(333, 17).Value = 25
$ws.Cells.Item(333, 18).Value = "Hortaliza"
